$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new employee row (row 4), copying the formatting/style of row 3 ---
$ws.Range("A3:W3").Copy($ws.Range("A4:W4"))

# Identity / classification columns (new shared strings, in the order they
# first appear: CAVIEDEZ FERNANDEZ, JUAN SEBASTIAN, OPERARIO, DIRECTA)
$ws.Range("A4").Value = "CAVIEDEZ FERNANDEZ"
$ws.Range("B4").Value = "JUAN SEBASTIAN"
$ws.Range("C4").Value = "OPERARIO"
$ws.Range("D4").Value = "DIRECTA"

# Numeric payroll columns for the new row
$ws.Range("E4").Value = 15333333.333333334
$ws.Range("F4").Value = 0.0
$ws.Range("G4").Value = 15333333.333333334
$ws.Range("H4").Value = 613333.3333333334
$ws.Range("I4").Value = 613333.3333333334
$ws.Range("J4").Value = 1226666.6666666667
$ws.Range("K4").Value = 1303333.3333333335
$ws.Range("L4").Value = 1840000.0
$ws.Range("M4").Value = 373520.0
$ws.Range("N4").Value = 3516853.3333333335
$ws.Range("O4").Value = 1277266.6666666667
$ws.Range("P4").Value = 1277266.6666666667
$ws.Range("Q4").Value = 639400.0
$ws.Range("R4").Value = 153333.33333333334
$ws.Range("S4").Value = 3347266.666666667
$ws.Range("T4").Value = 613333.3333333334
$ws.Range("U4").Value = 306666.6666666667
$ws.Range("V4").Value = 460000.0
$ws.Range("W4").Value = 1380000.0

# --- Recalculated ATEP/TOTAL (SEGURIDAD SOCIAL) figures on row 3 now that a
# second employee shares the accident-risk distribution ---
$ws.Range("M3").Value = 5202.599999999999
$ws.Range("N3").Value = 124802.59999999999

# --- Column widths widen (best fit) to accommodate the new, longer values ---
$ws.Columns.Item(1).ColumnWidth = 20.333333333333336
$ws.Columns.Item(2).ColumnWidth = 15.393229166666666
$ws.Columns.Item(5).ColumnWidth = 11.338541666666666
$ws.Columns.Item(7).ColumnWidth = 11.338541666666666
$ws.Columns.Item(8).ColumnWidth = 8.557291666666666
$ws.Columns.Item(9).ColumnWidth = 8.557291666666666
$ws.Columns.Item(10).ColumnWidth = 10.221354166666666
$ws.Columns.Item(11).ColumnWidth = 10.221354166666666
$ws.Columns.Item(12).ColumnWidth = 10.221354166666666
$ws.Columns.Item(13).ColumnWidth = 8.557291666666666
$ws.Columns.Item(14).ColumnWidth = 10.221354166666666
$ws.Columns.Item(19).ColumnWidth = 10.221354166666666
$ws.Columns.Item(21).ColumnWidth = 8.557291666666666
$ws.Columns.Item(22).ColumnWidth = 8.557291666666666
$ws.Columns.Item(23).ColumnWidth = 10.221354166666666
